# Update cryptos list with latest prices and 1h volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    # Force the cell to store $val as text (not auto-converted to a number),
    # then restore the default ("Normal") style so no stray formatting is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '43.534.65'
$ws.Range("E2").Value = '  +2.58%  '
Set-TextValue $ws.Range("D3") '2.417.95'
$ws.Range("E3").Value = '  +8.92%  '
$ws.Range("E4").Value = '  -0.32%  '
Set-TextValue $ws.Range("D5") '323.54'
$ws.Range("E5").Value = '  +13.01%  '
Set-TextValue $ws.Range("D6") '104.22'
$ws.Range("E6").Value = '  -4.92%  '
Set-TextValue $ws.Range("D7") '0.650'
$ws.Range("E7").Value = '  +4.78%  '
Set-TextValue $ws.Range("D9") '0.656'
$ws.Range("E9").Value = '  +10.32%  '
Set-TextValue $ws.Range("D10") '41.86'
$ws.Range("E10").Value = '  -3.08%  '
Set-TextValue $ws.Range("D11") '0.0949'
$ws.Range("E11").Value = '  +4.74%  '
Set-TextValue $ws.Range("D12") '8.70'
$ws.Range("E12").Value = '  +1.45%  '
Set-TextValue $ws.Range("D13") '1.04'
$ws.Range("E13").Value = '  +2.66%  '
Set-TextValue $ws.Range("D14") '17.44'
$ws.Range("E14").Value = '  +17.89%  '
$ws.Range("E15").Value = '  +2.83%  '
Set-TextValue $ws.Range("D16") '2.782.95'
$ws.Range("E16").Value = '  +9.01%  '
Set-TextValue $ws.Range("D17") '2.499.43'
$ws.Range("E17").Value = '  +11.76%  '
Set-TextValue $ws.Range("D18") '43.568.90'
$ws.Range("E18").Value = '  +3.02%  '
$ws.Range("E19").Value = '  +6.07%  '
Set-TextValue $ws.Range("D20") '7.45'
$ws.Range("E20").Value = '  +4.45%  '
Set-TextValue $ws.Range("D21") '75.64'
$ws.Range("E21").Value = '  +3.99%  '
Set-TextValue $ws.Range("D22") '3.49'
$ws.Range("E22").Value = '  +4.61%  '
Set-TextValue $ws.Range("D23") '260.78'
$ws.Range("E23").Value = '  +14.06%  '
Set-TextValue $ws.Range("D24") '2.45'
$ws.Range("E24").Value = '  +2.40%  '
Set-TextValue $ws.Range("D25") '9.67'
$ws.Range("E25").Value = '  +8.71%  '
Set-TextValue $ws.Range("D26") '12.01'
$ws.Range("E26").Value = '  +5.65%  '
Set-TextValue $ws.Range("D27") '0.999'
$ws.Range("E27").Value = '  +0.01%  '
Set-TextValue $ws.Range("D28") '22.97'
$ws.Range("E28").Value = '  +10.62%  '
$ws.Range("E29").Value = '  +2.11%  '
Set-TextValue $ws.Range("D30") '178.24'
$ws.Range("E30").Value = '  +3.33%  '
Set-TextValue $ws.Range("D31") '38.19'
$ws.Range("E31").Value = '  +4.12%  '
Set-TextValue $ws.Range("D32") '3.24'
$ws.Range("E32").Value = '  +1.60%  '
Set-TextValue $ws.Range("D33") '0.0937'
$ws.Range("E33").Value = '  +7.77%  '
$ws.Range("E34").Value = '  +7.76%  '
$ws.Range("E35").Value = '  +6.35%  '
Set-TextValue $ws.Range("D36") '4.88'
$ws.Range("E36").Value = '  -1.29%  '
Set-TextValue $ws.Range("D37") '0.0372'
$ws.Range("E37").Value = '  +0.93%  '
Set-TextValue $ws.Range("D38") '3.94'
$ws.Range("E38").Value = '  -4.88%  '
Set-TextValue $ws.Range("D41") '1.66'
$ws.Range("E41").Value = '  +28.93%  '
Set-TextValue $ws.Range("D42") '0.235'
$ws.Range("E42").Value = '  +3.27%  '
Set-TextValue $ws.Range("D43") '124.81'
$ws.Range("E43").Value = '  +23.94%  '
Set-TextValue $ws.Range("D44") '69.41'
$ws.Range("E44").Value = '  -6.64%  '
$ws.Range("E45").Value = '  +0.02%  '
Set-TextValue $ws.Range("D46") '12.65'
$ws.Range("E46").Value = '  +3.44%  '
$ws.Range("E49").Value = '  +4.12%  '
Set-TextValue $ws.Range("D50") '1.601.56'
$ws.Range("E50").Value = '  +14.22%  '
$ws.Range("E51").Value = '  +3.73%  '

# Row 39/40 swap: Kaspa <-> LidoDAOToken with updated values
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D39") "2.92"
$ws.Range("E39").Value = "  +23.04%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D40") "0.106"
$ws.Range("E40").Value = "  +1.87%  "

# Row 47/48 swap: FraxShare <-> THORChain with updated values
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D47") "5.70"
$ws.Range("E47").Value = "  +6.36%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D48") "9.54"
$ws.Range("E48").Value = "  +13.65%  "
